# Insert two new data rows before row 25 (shifts existing rows 25-124 down to 27-126)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A25:T26").Insert()

# New row 25: Kurakata / Primera
$ws.Cells.Item(25, 1).Value = 11
$ws.Cells.Item(25, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(25, 3).Value = "Bíobío"
$ws.Cells.Item(25, 4).Value = 44560
$ws.Cells.Item(25, 5).Value = 8
$ws.Cells.Item(25, 6).Value = "Fruta"
$ws.Cells.Item(25, 7).Value = 100103
$ws.Cells.Item(25, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(25, 9).Value = 100103004
$ws.Cells.Item(25, 10).Value = "Durazno"
$ws.Cells.Item(25, 11).Value = "Kurakata"
$ws.Cells.Item(25, 12).Value = "Primera"
$ws.Cells.Item(25, 13).Value = 100
$ws.Cells.Item(25, 14).Value = 10000
$ws.Cells.Item(25, 15).Value = 11000
$ws.Cells.Item(25, 16).Value = 10500
$ws.Cells.Item(25, 17).Value = "$/caja 16 kilos empedrada"
$ws.Cells.Item(25, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(25, 19).Value = 656
$ws.Cells.Item(25, 20).Value = 16

# New row 26: Kurakata / Segunda
$ws.Cells.Item(26, 1).Value = 11
$ws.Cells.Item(26, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(26, 3).Value = "Bíobío"
$ws.Cells.Item(26, 4).Value = 44560
$ws.Cells.Item(26, 5).Value = 8
$ws.Cells.Item(26, 6).Value = "Fruta"
$ws.Cells.Item(26, 7).Value = 100103
$ws.Cells.Item(26, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(26, 9).Value = 100103004
$ws.Cells.Item(26, 10).Value = "Durazno"
$ws.Cells.Item(26, 11).Value = "Kurakata"
$ws.Cells.Item(26, 12).Value = "Segunda"
$ws.Cells.Item(26, 13).Value = 50
$ws.Cells.Item(26, 14).Value = 9000
$ws.Cells.Item(26, 15).Value = 9000
$ws.Cells.Item(26, 16).Value = 9000
$ws.Cells.Item(26, 17).Value = "$/caja 16 kilos empedrada"
$ws.Cells.Item(26, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(26, 19).Value = 562
$ws.Cells.Item(26, 20).Value = 16
